# Update the column headers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "building_name"
$ws.Range("B1").Value = "building_no"
$ws.Range("C1").Value = "website:map"

# Clean up the URL values in column C (rows 2-15): remove the stray space
# after "sharepoi=" and any trailing whitespace.
for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cleaned = $current.Trim()
        $cleaned = $cleaned.Replace("sharepoi= ", "sharepoi=")
        $cell.Value = $cleaned
    }
}
